# Auto-generated update of Leve price/profit snapshot columns (H:N)
# across the Diabolos_Profits workbook, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 53154.7
$ws.Range("I28").Value = 59923.707
$ws.Range("J28").Value = 14797
$ws.Range("K28").Value = 59923.707
$ws.Range("L28").Value = 14797
$ws.Range("M28").Value = -59438.707
$ws.Range("N28").Value = -15767
$ws.Range("H29").Value = 6500
$ws.Range("J29").Value = 6500
$ws.Range("L29").Value = 19500
$ws.Range("N29").Value = -20062
$ws.Range("H33").Value = 45226.53
$ws.Range("I33").Value = 58907.46
$ws.Range("K33").Value = 58907.46
$ws.Range("M33").Value = -58678.46
$ws.Range("H38").Value = 153.5
$ws.Range("I38").Value = 153.5
$ws.Range("K38").Value = 460.5
$ws.Range("M38").Value = -88.5
$ws.Range("H58").Value = 1155.4
$ws.Range("I58").Value = 364.85715
$ws.Range("K58").Value = 1094.57145
$ws.Range("M58").Value = -944.5714499999999
$ws.Range("H138").Value = 2781.2903
$ws.Range("I138").Value = 2100.2856
$ws.Range("J138").Value = 3342.1177
$ws.Range("K138").Value = 6300.8568
$ws.Range("L138").Value = 10026.3531
$ws.Range("M138").Value = -1160.8568
$ws.Range("N138").Value = -20306.3531

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 7900
$ws.Range("I6").Value = 800
$ws.Range("K6").Value = 800
$ws.Range("M6").Value = -627
$ws.Range("H122").Value = 4483.857
$ws.Range("I122").Value = 4019.4546
$ws.Range("K122").Value = 12058.3638
$ws.Range("M122").Value = -9608.363799999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1938.6666
$ws.Range("I105").Value = 1853.6471
$ws.Range("K105").Value = 1853.6471
$ws.Range("M105").Value = -106.6470999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 68515.46000000001
$ws.Range("I86").Value = 95997.78
$ws.Range("K86").Value = 95997.78
$ws.Range("M86").Value = -94874.78
$ws.Range("H89").Value = 68515.46000000001
$ws.Range("I89").Value = 95997.78
$ws.Range("K89").Value = 479988.9
$ws.Range("M89").Value = -474372.9
$ws.Range("H132").Value = 2539.75
$ws.Range("I132").Value = 1744.1
$ws.Range("J132").Value = 4528.875
$ws.Range("K132").Value = 5232.299999999999
$ws.Range("L132").Value = 13586.625
$ws.Range("M132").Value = -2702.299999999999
$ws.Range("N132").Value = -18646.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1401.9474
$ws.Range("I23").Value = 925.1111
$ws.Range("K23").Value = 2775.3333
$ws.Range("M23").Value = -2540.3333
$ws.Range("H28").Value = 2988.3333
$ws.Range("I28").Value = 2988.3333
$ws.Range("K28").Value = 8964.999899999999
$ws.Range("M28").Value = -8732.999899999999
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("K94").Value = 6000
$ws.Range("M94").Value = -5324
$ws.Range("H98").Value = 4798.0586
$ws.Range("I98").Value = 9049.286
$ws.Range("J98").Value = 1822.2
$ws.Range("K98").Value = 27147.858
$ws.Range("L98").Value = 5466.6
$ws.Range("M98").Value = -25649.858
$ws.Range("N98").Value = -8462.6
$ws.Range("H103").Value = 1431.5714
$ws.Range("I103").Value = 1025
$ws.Range("J103").Value = 1499.3334
$ws.Range("K103").Value = 3075
$ws.Range("L103").Value = 4498.0002
$ws.Range("M103").Value = -2196
$ws.Range("N103").Value = -6256.0002
$ws.Range("H118").Value = 1250.4
$ws.Range("I118").Value = 1250.4
$ws.Range("K118").Value = 3751.2
$ws.Range("M118").Value = -2508.2
$ws.Range("H139").Value = 1758.8182
$ws.Range("I139").Value = 1049.5714
$ws.Range("K139").Value = 3148.7142
$ws.Range("M139").Value = 1991.2858
$ws.Range("H140").Value = 3004.8
$ws.Range("I140").Value = 2512
$ws.Range("K140").Value = 7536
$ws.Range("M140").Value = -2356

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3286.1428
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 3200.5
$ws.Range("K80").Value = 3800
$ws.Range("L80").Value = 3200.5
$ws.Range("M80").Value = -2802
$ws.Range("N80").Value = -5196.5
$ws.Range("H83").Value = 3286.1428
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 3200.5
$ws.Range("K83").Value = 19000
$ws.Range("L83").Value = 16002.5
$ws.Range("M83").Value = -14008
$ws.Range("N83").Value = -25986.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23812796
$ws.Range("I7").Value = 38464150
$ws.Range("K7").Value = 38464150
$ws.Range("M7").Value = -38464038
$ws.Range("H22").Value = 6456212.5
$ws.Range("J22").Value = 8069266
$ws.Range("L22").Value = 8069266
$ws.Range("N22").Value = -8069856
$ws.Range("H27").Value = 6456212.5
$ws.Range("J27").Value = 8069266
$ws.Range("L27").Value = 8069266
$ws.Range("N27").Value = -8069480
$ws.Range("H55").Value = 407.3158
$ws.Range("I55").Value = 400.125
$ws.Range("J55").Value = 445.66666
$ws.Range("K55").Value = 400.125
$ws.Range("L55").Value = 445.66666
$ws.Range("M55").Value = -227.125
$ws.Range("N55").Value = -791.66666
$ws.Range("H82").Value = 1397.5714
$ws.Range("I82").Value = 1026.6666
$ws.Range("J82").Value = 1675.75
$ws.Range("K82").Value = 1026.6666
$ws.Range("L82").Value = 1675.75
$ws.Range("M82").Value = -665.6666
$ws.Range("N82").Value = -2397.75
$ws.Range("H85").Value = 1397.5714
$ws.Range("I85").Value = 1026.6666
$ws.Range("J85").Value = 1675.75
$ws.Range("K85").Value = 1026.6666
$ws.Range("L85").Value = 1675.75
$ws.Range("M85").Value = 221.3334
$ws.Range("N85").Value = -4171.75
$ws.Range("H126").Value = 23812796
$ws.Range("I126").Value = 38464150
$ws.Range("K126").Value = 115392450
$ws.Range("M126").Value = -115389980
$ws.Range("H132").Value = 6257.8335
$ws.Range("I132").Value = 3437.6924
$ws.Range("K132").Value = 10313.0772
$ws.Range("M132").Value = -7783.0772
$ws.Range("H136").Value = 6432.6665
$ws.Range("I136").Value = 1762.625
$ws.Range("J136").Value = 11769.857
$ws.Range("K136").Value = 5287.875
$ws.Range("L136").Value = 35309.571
$ws.Range("M136").Value = -2737.875
$ws.Range("N136").Value = -40409.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = -4860
$ws.Range("N13").Value = -5280
$ws.Range("H100").Value = 465.33334
$ws.Range("I100").Value = 465.33334
$ws.Range("K100").Value = 930.66668
$ws.Range("M100").Value = -389.66668
$ws.Range("H132").Value = 265782.75
$ws.Range("I132").Value = 315042.06
$ws.Range("J132").Value = 3066.5
$ws.Range("K132").Value = 945126.1799999999
$ws.Range("L132").Value = 9199.5
$ws.Range("M132").Value = -942596.1799999999
$ws.Range("N132").Value = -14259.5
$ws.Range("H136").Value = 4923.2085
$ws.Range("J136").Value = 7960.857
$ws.Range("L136").Value = 23882.571
$ws.Range("N136").Value = -28982.571
